$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Utility (D column) values for the rows that changed
$ws.Range("D2").Value = 0.46125965385188411
$ws.Range("D3").Value = 0.17650348917363121
$ws.Range("D4").Value = 0.16586840448085979
$ws.Range("D5").Value = 0.15812655143109419
$ws.Range("D6").Value = 0.13690883275877391
$ws.Range("D7").Value = 0.1323715121629416
$ws.Range("D8").Value = 0.1312068865157234
$ws.Range("D12").Value = 0.11747580873334899
$ws.Range("D14").Value = 0.1169297340348983
$ws.Range("D15").Value = 0.11546990160573629
$ws.Range("D16").Value = 0.11294464072590001
$ws.Range("D17").Value = 0.1102228301242768
$ws.Range("D18").Value = 0.1083854215234004
$ws.Range("D19").Value = 0.086160824674101674
$ws.Range("D20").Value = 0.084872465027004032
$ws.Range("D21").Value = 0.077638588852357124
$ws.Range("D22").Value = 0.072098318752407742
$ws.Range("D23").Value = 0.066935055885867442
$ws.Range("D24").Value = 0.042267933562189643
$ws.Range("D26").Value = 0.01551863742923161
$ws.Range("D27").Value = 0.0125224602008272

# Rows 17 and 18 also swapped their property_type (A) and Function (C) values
$ws.Range("A17").Value = "property_type"
$ws.Range("C17").Value = "avg"
$ws.Range("A18").Value = "zipcode"
$ws.Range("C18").Value = "sum"

# Update the selected range shown in the workbook
$ws.Range("A2:XFD10").Select()
